# Updated symbol list on Sun Dec 18 16:52:24 UTC 2022 with GitHub Actions
#
# Applies the cryptos.xlsx "Price" (column D) refresh, the ProBitToken
# "Bestin24h" tag addition, and the CEJI / BKEXToken row-content swap
# described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    # Values that look numeric ("245.52", "0.05618", ...) need a leading
    # quote so Excel stores them as text (matching the source file's
    # inline-string cells) instead of silently re-typing them as numbers.
    $needsQuote = $Text -match '^-?[0-9]*\.?[0-9]+([eE][-+]?[0-9]+)?$'
    if ($needsQuote) {
        $cell.Value = "'" + $Text
    } else {
        $cell.Value = $Text
    }
}

# ---- Price (column D) updates ----
Set-TextValue "D2"  "245.52"
Set-TextValue "D4"  "5.469"
Set-TextValue "D5"  "0.05618"
Set-TextValue "D7"  "6.464"
Set-TextValue "D8"  "0.8049"
Set-TextValue "D9"  "1.042"
Set-TextValue "D10" "0.1421"
Set-TextValue "D11" "0.07318"
Set-TextValue "D12" "0.03182"
Set-TextValue "D13" "0.02940"
Set-TextValue "D14" "0.09272"
Set-TextValue "D15" "0.001663"
Set-TextValue "D16" "3.216"
Set-TextValue "D17" "0.04737"
Set-TextValue "D18" "0.0005827"
Set-TextValue "D19" "0.006457"
Set-TextValue "D20" "0.005068"
Set-TextValue "D22" "0.0001504"
Set-TextValue "D23" "3.989"
Set-TextValue "D25" "0.3293"

# ---- ProBitToken (row 26) ----
Set-TextValue "D26" "0.1254"
Set-TextValue "E26" "25ProBitTokenPROBBestin24h"

Set-TextValue "D27" "0.0002909"
Set-TextValue "D40" "0.04157"
Set-TextValue "D41" "0.006886"

# ---- Row 42 / 43: CEJI and BKEXToken swap places ----
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1041"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002977"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D44" "0.009541"
Set-TextValue "D45" "0.00005658"
Set-TextValue "D47" "0.6820"
Set-TextValue "D48" "0.01451"
Set-TextValue "D49" "0.00002106"
Set-TextValue "D50" "0.01013"
